$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Update progress/percentage values for iteration 6 tasks (part of the 90%)
$ws.Range("G11").Value = 90
$ws.Range("G12").Value = 90
$ws.Range("G23").Value = 90
$ws.Range("G31").Value = 90

# Reposition the view: scroll so row 12 is at the top, and select G31
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G31").Select()
